$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title block updates ---
$ws.Range("A8").Value = "Volume 30   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/13/2023  Through  11/19/2023"

# --- Crime complaint table (rows 14-30): new week's figures ---
$ws.Range("A14").Value = "Murder"
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = -50
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = -66.666666666666
$ws.Range("I14").Value = 18
$ws.Range("J14").Value = 26
$ws.Range("K14").Value = -30.769230769230
$ws.Range("L14").Value = -30.769230769230
$ws.Range("M14").Value = -35.714285714285
$ws.Range("N14").Value = -83.333333333333
$ws.Range("A15").Value = "Rape"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 75
$ws.Range("F15").Value = 14
$ws.Range("G15").Value = 24
$ws.Range("H15").Value = -41.666666666666
$ws.Range("I15").Value = 178
$ws.Range("J15").Value = 187
$ws.Range("K15").Value = -4.812834224598
$ws.Range("L15").Value = 21.088435374149
$ws.Range("M15").Value = 47.107438016528
$ws.Range("N15").Value = -3.260869565217
$ws.Range("A16").Value = "Robbery"
$ws.Range("C16").Value = 35
$ws.Range("D16").Value = 46
$ws.Range("E16").Value = -23.913043478260
$ws.Range("F16").Value = 170
$ws.Range("G16").Value = 197
$ws.Range("H16").Value = -13.705583756345
$ws.Range("I16").Value = 1824
$ws.Range("J16").Value = 1664
$ws.Range("K16").Value = 9.615384615384
$ws.Range("L16").Value = 73.055028462998
$ws.Range("M16").Value = 7.231040564373
$ws.Range("N16").Value = -76.657281801894
$ws.Range("A17").Value = "Fel. Assault"
$ws.Range("C17").Value = 54
$ws.Range("D17").Value = 48
$ws.Range("E17").Value = 12.5
$ws.Range("F17").Value = 233
$ws.Range("G17").Value = 207
$ws.Range("H17").Value = 12.56038647343
$ws.Range("I17").Value = 2597
$ws.Range("J17").Value = 2232
$ws.Range("K17").Value = 16.353046594982
$ws.Range("L17").Value = 45.083798882681
$ws.Range("M17").Value = 90.535583272193
$ws.Range("N17").Value = -0.192159877017
$ws.Range("A18").Value = "Burglary"
$ws.Range("C18").Value = 36
$ws.Range("D18").Value = 46
$ws.Range("E18").Value = -21.739130434782
$ws.Range("F18").Value = 149
$ws.Range("G18").Value = 170
$ws.Range("H18").Value = -12.352941176470
$ws.Range("I18").Value = 1735
$ws.Range("J18").Value = 1767
$ws.Range("K18").Value = -1.810979060554
$ws.Range("L18").Value = 18.673050615595
$ws.Range("M18").Value = -25.600343053173
$ws.Range("N18").Value = -86.866010598031
$ws.Range("A19").Value = "Gr. Larceny"
$ws.Range("C19").Value = 124
$ws.Range("D19").Value = 117
$ws.Range("E19").Value = 5.982905982905
$ws.Range("F19").Value = 495
$ws.Range("G19").Value = 533
$ws.Range("H19").Value = -7.129455909943
$ws.Range("I19").Value = 5988
$ws.Range("J19").Value = 6281
$ws.Range("K19").Value = -4.664862283075
$ws.Range("L19").Value = 54.968944099378
$ws.Range("M19").Value = 66.889632107023
$ws.Range("N19").Value = -19.277433270423
$ws.Range("A20").Value = "G.L.A."
$ws.Range("C20").Value = 58
$ws.Range("D20").Value = 40
$ws.Range("E20").Value = 45
$ws.Range("F20").Value = 191
$ws.Range("G20").Value = 178
$ws.Range("H20").Value = 7.303370786516
$ws.Range("I20").Value = 2398
$ws.Range("J20").Value = 1718
$ws.Range("K20").Value = 39.580908032596
$ws.Range("L20").Value = 94.48499594485
$ws.Range("M20").Value = 52.641629535327
$ws.Range("N20").Value = -88.640992847330
$ws.Range("A21").Value = "TOTAL"
$ws.Range("C21").Value = 315
$ws.Range("D21").Value = 303
$ws.Range("E21").Value = 3.960396039603
$ws.Range("F21").Value = 1253
$ws.Range("G21").Value = 1312
$ws.Range("H21").Value = -4.496951219512
$ws.Range("I21").Value = 14738
$ws.Range("J21").Value = 13875
$ws.Range("K21").Value = 6.219819819819
$ws.Range("L21").Value = 53.90559732665
$ws.Range("M21").Value = 37.686846038864
$ws.Range("N21").Value = -71.899250672107
$ws.Range("A22").Value = "Transit"
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 28
$ws.Range("G22").Value = 28
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 322
$ws.Range("J22").Value = 276
$ws.Range("K22").Value = 16.666666666666
$ws.Range("L22").Value = 128.368794326241
$ws.Range("M22").Value = 86.127167630057
$ws.Range("N22").Value = "***.*"
$ws.Range("A23").Value = "Housing"
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 19
$ws.Range("G23").Value = 25
$ws.Range("H23").Value = -24
$ws.Range("I23").Value = 241
$ws.Range("J23").Value = 216
$ws.Range("K23").Value = 11.574074074074
$ws.Range("L23").Value = 20.5
$ws.Range("M23").Value = 57.516339869281
$ws.Range("N23").Value = "***.*"
$ws.Range("A24").Value = "Petit Larceny"
$ws.Range("C24").Value = 346
$ws.Range("D24").Value = 371
$ws.Range("E24").Value = -6.738544474393
$ws.Range("F24").Value = 1164
$ws.Range("G24").Value = 1257
$ws.Range("H24").Value = -7.398568019093
$ws.Range("I24").Value = 13739
$ws.Range("J24").Value = 13619
$ws.Range("K24").Value = 0.881121961964
$ws.Range("L24").Value = 27.995155580398
$ws.Range("M24").Value = 62.938804554079
$ws.Range("N24").Value = "***.*"
$ws.Range("A25").Value = "Misd. Assault"
$ws.Range("C25").Value = 124
$ws.Range("D25").Value = 95
$ws.Range("E25").Value = 30.526315789473
$ws.Range("F25").Value = 465
$ws.Range("G25").Value = 425
$ws.Range("H25").Value = 9.411764705882
$ws.Range("I25").Value = 4807
$ws.Range("J25").Value = 4444
$ws.Range("K25").Value = 8.168316831683
$ws.Range("L25").Value = 23.256410256410
$ws.Range("M25").Value = 13.026099224077
$ws.Range("N25").Value = "***.*"
$ws.Range("A26").Value = "UCR Rape*"
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 80
$ws.Range("F26").Value = 21
$ws.Range("G26").Value = 32
$ws.Range("H26").Value = -34.375
$ws.Range("I26").Value = 279
$ws.Range("J26").Value = 267
$ws.Range("K26").Value = 4.494382022471
$ws.Range("L26").Value = 14.344262295082
$ws.Range("M26").Value = "***.*"
$ws.Range("N26").Value = "***.*"
$ws.Range("A27").Value = "Other Sex Crimes"
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 14
$ws.Range("E27").Value = -57.142857142857
$ws.Range("F27").Value = 48
$ws.Range("G27").Value = 63
$ws.Range("H27").Value = -23.809523809523
$ws.Range("I27").Value = 627
$ws.Range("J27").Value = 582
$ws.Range("K27").Value = 7.731958762886
$ws.Range("L27").Value = 26.156941649899
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"
$ws.Range("A28").Value = "Shooting Vic."
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -66.666666666666
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -80
$ws.Range("I28").Value = 60
$ws.Range("J28").Value = 68
$ws.Range("K28").Value = -11.764705882352
$ws.Range("L28").Value = -20
$ws.Range("M28").Value = 39.534883720930
$ws.Range("N28").Value = -74.137931034482
$ws.Range("A29").Value = "Shooting Inc."
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = -66.666666666666
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = -80
$ws.Range("I29").Value = 55
$ws.Range("J29").Value = 58
$ws.Range("K29").Value = -5.172413793103
$ws.Range("L29").Value = 3.773584905660
$ws.Range("M29").Value = 57.142857142857
$ws.Range("N29").Value = -73.429951690821
$ws.Range("A30").Value = "Hate Crimes"
# C30 holds the literal text "0" (not the number 0) in the source data,
# mirroring how C14/C28/C29 used to store it as text before this week's update.
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = 9
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = 200
$ws.Range("I30").Value = 66
$ws.Range("J30").Value = 57
$ws.Range("K30").Value = 15.789473684210
$ws.Range("L30").Value = -1.492537313432
$ws.Range("M30").Value = "***.*"
$ws.Range("N30").Value = "***.*"

# --- Historical Perspective section (rows 37-43): labels only shift up one
#     slot because the "Murder" row no longer appears here; values unchanged ---
$ws.Range("A37").Value = "Rape"
$ws.Range("A38").Value = "Robbery"
$ws.Range("A39").Value = "Fel. Assault"
$ws.Range("A40").Value = "Burglary"
$ws.Range("A41").Value = "Gr. Larceny"
$ws.Range("A42").Value = "G.L.A."
$ws.Range("A43").Value = "TOTAL"
